$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume/number and week-covering date range) ---
$ws.Range("A8").Value = "Volume 30   Number  42"
$ws.Range("C9").Value = "Report Covering the Week  10/16/2023  Through  10/22/2023"

# --- Plain numeric value updates (same cell type, new number) ---
$numericUpdates = @{
    "F14" = 1
    "L14" = 71.428571428571
    "M15" = -50
    "C16" = 3
    "D16" = 2
    "E16" = 50
    "F16" = 9
    "G16" = 7
    "H16" = 28.571428571428
    "I16" = 126
    "J16" = 88
    "K16" = 43.181818181818
    "L16" = 80
    "M16" = -46.382978723404
    "N16" = -84.558823529411
    "C17" = 8
    "D17" = 12
    "E17" = -33.333333333333
    "F17" = 31
    "G17" = 31
    "H17" = 0
    "I17" = 378
    "J17" = 324
    "K17" = 16.666666666666
    "L17" = 66.519823788546
    "M17" = 52.419354838709
    "N17" = -42.024539877300
    "C18" = 3
    "E18" = 50
    "F18" = 8
    "G18" = 7
    "H18" = 14.285714285714
    "I18" = 102
    "J18" = 90
    "K18" = 13.333333333333
    "L18" = 100
    "M18" = -59.362549800796
    "N18" = -92.816901408450
    "F19" = 34
    "G19" = 20
    "H19" = 70
    "I19" = 340
    "J19" = 277
    "K19" = 22.743682310469
    "L19" = 55.251141552511
    "M19" = 2.409638554216
    "N19" = -26.406926406926
    "F20" = 16
    "G20" = 13
    "H20" = 23.076923076923
    "I20" = 129
    "K20" = 37.234042553191
    "L20" = 111.475409836066
    "M20" = -21.341463414634
    "N20" = -88.022284122562
    "C21" = 25
    "E21" = 13.636363636363
    "F21" = 99
    "G21" = 79
    "H21" = 25.316455696202
    "I21" = 1100
    "J21" = 890
    "K21" = 23.595505617977
    "L21" = 68.970814132104
    "M21" = -13.180741910023
    "N21" = -75.495656048117
    "E23" = -100
    "F23" = 5
    "G23" = 3
    "H23" = 66.666666666666
    "J23" = 48
    "K23" = 45.833333333333
    "L23" = 70.731707317073
    "C24" = 17
    "D24" = 25
    "E24" = -32
    "F24" = 94
    "H24" = -5.050505050505
    "I24" = 1012
    "J24" = 965
    "K24" = 4.870466321243
    "L24" = 55.692307692307
    "M24" = -18.714859437751
    "C25" = 21
    "D25" = 7
    "E25" = 200
    "F25" = 67
    "G25" = 41
    "H25" = 63.414634146341
    "I25" = 595
    "J25" = 523
    "K25" = 13.766730401529
    "L25" = 33.707865168539
    "M25" = -42.788461538461
    "G27" = 4
    "H27" = 25
    "L27" = 72.093023255813
    "G28" = 1
    "H28" = 0
    "L28" = -20
    "N28" = -77.011494252873
    "G29" = 1
    "H29" = 0
    "L29" = -21.739130434782
    "N29" = -76
    "F30" = 1
}
foreach ($ref in $numericUpdates.Keys) {
    $ws.Range($ref).Value = $numericUpdates[$ref]
}

# --- Cells that change from a plain number to the special text markers
#     "0" (no data / zero placeholder) or "***.*" (change undefined).
#     A14 already carries the target style (s="14"), so its format is
#     copied onto each converted cell to keep formatting identical; the
#     leading apostrophe forces Excel to store the value as text instead
#     of re-parsing "0" back into a number. ---
$formatDonor = $ws.Range("A14")
$textUpdates = @{
    "G14" = "0"
    "H14" = "***.*"
    "D15" = "0"
    "E15" = "***.*"
    "D20" = "0"
    "E20" = "***.*"
    "C23" = "0"
    "C26" = "0"
    "D26" = "0"
    "E26" = "***.*"
    "C27" = "0"
    "D27" = "0"
    "E27" = "***.*"
    "G30" = "0"
    "H30" = "***.*"
}
foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = "'" + $textUpdates[$ref]
    $formatDonor.Copy()
    $ws.Range($ref).PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}
